# "added 4wk low sales check"
# Update Inventory Coverage (H) and Seasonality Index (L) values on the
# "Forecast Comparison" sheet, and the derived Total Forecast summary
# figures (16-week and 4-week) on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# row -> (Inventory Coverage, Seasonality Index)
$updates = @(
    @{ Row = 2;  H = 30.71; L = 1.17 },
    @{ Row = 3;  H = 26;    L = 0.9  },
    @{ Row = 4;  H = 25;    L = 0.92 },
    @{ Row = 5;  H = 27.43; L = 0.88 },
    @{ Row = 6;  H = 30.83; L = 1.01 },
    @{ Row = 7;  H = 29.83; L = 0.93 },
    @{ Row = 8;  H = 24.71; L = 1.03 },
    @{ Row = 9;  H = 20.75; L = 1.17 },
    @{ Row = 10; H = 19.75; L = 1.1  },
    @{ Row = 11; H = 21.43; L = 0.99 },
    @{ Row = 12; H = 23.83; L = 1    },
    @{ Row = 13; H = 22.83; L = 1    },
    @{ Row = 14; H = 18.71; L = 1.19 },
    @{ Row = 15; H = 15.5  },
    @{ Row = 16; H = 16.57; L = 0.89 },
    @{ Row = 17; H = 15.57; L = 1.15 }
)

foreach ($u in $updates) {
    $wsForecast.Cells.Item($u.Row, 8).Value = $u.H
    if ($u.ContainsKey("L")) {
        $wsForecast.Cells.Item($u.Row, 12).Value = $u.L
    }
}

# Summary sheet totals (stored as text in the workbook)
$wsSummary.Range("B9").Value = "23"
$wsSummary.Range("B11").Value = "6"
